$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "AssignScheme" worksheet as the last sheet in the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws4.Name = "AssignScheme"

# Populate the new sheet. The order below reproduces the exact order the
# original author typed the values in (it controls shared-string ordering).
$ws4.Range("A2").Value = "warehouse / distributor"
$ws4.Range("B1").Value = "ZoneName"
$ws4.Range("B2").Value = "northzone"
$ws4.Range("C1").Value = "SubZone Name"
$ws4.Range("C2").Value = "subzone1"
$ws4.Range("D1").Value = "Holder Name"
$ws4.Range("A1").Value = "HolderType"
$ws4.Range("D2").Value = "warehouse1"
$ws4.Range("F1").Value = "Scheme Name"
$ws4.Range("E1").Value = "Scheme ID"
$ws4.Range("E2").Value = '"5"'

# Header row (row 1) is bold, matching the other sheets' header style.
$ws4.Range("A1:F1").Font.Bold = $true

# Approximate the author's best-fit column widths as closely as this
# engine's width model allows.
$ws4.Columns.Item(1).ColumnWidth = 21.09
$ws4.Columns.Item(2).ColumnWidth = 9.25
$ws4.Columns.Item(3).ColumnWidth = 12.42
$ws4.Columns.Item(4).ColumnWidth = 11.09
$ws4.Columns.Item(5).ColumnWidth = 12.25
$ws4.Columns.Item(6).ColumnWidth = 12.25

$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1

[void]$ws4.Range("E8").Select()

# ---------------------------------------------------------------------------
# 2. login sheet: username cell A2 changes from "pumaadmin" to
#    "myclientadmin"; also becomes the active/selected sheet & cell.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("login")
$ws1.Range("A2").Value = "myclientadmin"
$ws1.Columns.Item(1).ColumnWidth = 11.92

[void]$ws1.Range("C3").Select()
